# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 13:22"

# --- Reorder Portugal/Brasil (Portugal now has more cases than Brasil) ---
# Row 18 becomes Portugal, Row 19 becomes Brasil, with refreshed stats.
$ws.Range("A18").Value = "Portugal"
$ws.Range("B18").Value = 11730
$ws.Range("C18").Value = 452
$ws.Range("D18").Value = 140
$ws.Range("E18").Value = 11279
$ws.Range("F18").Value = 270
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 311

$ws.Range("A19").Value = "Brasil"
$ws.Range("B19").Value = 11281
$ws.Range("C19").Value = 27
$ws.Range("D19").Value = 127
$ws.Range("E19").Value = 10667
$ws.Range("F19").Value = 296
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 487

# --- Austria (row 17) refreshed stats ---
$ws.Range("B17").Value = 12109
$ws.Range("C17").Value = 58
$ws.Range("E17").Value = 8426

# --- Noruega (row 25) refreshed stats ---
$ws.Range("F25").Value = 83

# --- Dinamarca (row 27) refreshed stats ---
$ws.Range("D27").Value = 1378
$ws.Range("E27").Value = 3082
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 187

# --- Kuwait (row 74) refreshed stats ---
$ws.Range("F74").Value = 20

# --- Kazajistan (row 76) refreshed stats ---
$ws.Range("B76").Value = 629
$ws.Range("C76").Value = 45
$ws.Range("E76").Value = 578
